$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.152912314591731
$ws.Range("C2").Value = 0.09153488514983366
$ws.Range("D2").Value = 0.5257931236919404
$ws.Range("E2").Value = 0.1728898487075377
$ws.Range("G2").Value = 0.002565892654596835
$ws.Range("J2").Value = 0.06832942028598588
$ws.Range("K2").Value = 0.6084244706602533
$ws.Range("L2").Value = 0.3817415372010515
$ws.Range("M2").Value = 0.3334672989715379
$ws.Range("N2").Value = 3.172810236860332
$ws.Range("O2").Value = 7.603997695886903
# Row 3
$ws.Range("B3").Value = 1.122607785687222
$ws.Range("C3").Value = 0.09014219837622051
$ws.Range("D3").Value = 0.5244424576332847
$ws.Range("E3").Value = 0.1734020414272202
$ws.Range("G3").Value = 0.002568760308303305
$ws.Range("J3").Value = 0.06837021560265555
$ws.Range("K3").Value = 0.579298000960506
$ws.Range("L3").Value = 0.3802693966747839
$ws.Range("M3").Value = 0.3281998855934134
$ws.Range("N3").Value = 3.196780007223172
$ws.Range("O3").Value = 7.632383548754603
# Row 4
$ws.Range("B4").Value = 1.104511427369232
$ws.Range("C4").Value = 0.0892737462161719
$ws.Range("D4").Value = 0.5238164316904204
$ws.Range("E4").Value = 0.1737648938448082
$ws.Range("G4").Value = 0.002570616732927786
$ws.Range("J4").Value = 0.06839728786289756
$ws.Range("K4").Value = 0.561676104069349
$ws.Range("L4").Value = 0.3795044906745488
$ws.Range("M4").Value = 0.3251069130012034
$ws.Range("N4").Value = 3.21228313511385
$ws.Range("O4").Value = 7.652638335482976
# Row 5
$ws.Range("B5").Value = 1.097266079878295
$ws.Range("C5").Value = 0.08891649674404789
$ws.Range("D5").Value = 0.5236125631315645
$ws.Range("E5").Value = 0.1739249486411776
$ws.Range("G5").Value = 0.002571397370575921
$ws.Range("J5").Value = 0.06840883039371937
$ws.Range("K5").Value = 0.554561339224577
$ws.Range("L5").Value = 0.3792278388239154
$ws.Range("M5").Value = 0.3238821601453985
$ws.Range("N5").Value = 3.21879843538596
$ws.Range("O5").Value = 7.661603646261881
# Row 6
$ws.Range("B6").Value = 1.096070809371554
$ws.Range("C6").Value = 0.08885697368295098
$ws.Range("D6").Value = 0.5235818094470659
$ws.Range("E6").Value = 0.1739522626564352
$ws.Range("G6").Value = 0.002571528454254532
$ws.Range("J6").Value = 0.0684107778904739
$ws.Range("K6").Value = 0.5533839556637759
$ws.Range("L6").Value = 0.3791840211566679
$ws.Range("M6").Value = 0.3236809487551113
$ws.Range("N6").Value = 3.219892236759947
$ws.Range("O6").Value = 7.663135313768237
# Row 7
$ws.Range("B7").Value = 1.104413190656402
$ws.Range("C7").Value = 0.08926894177336209
$ws.Range("D7").Value = 0.5238134746055891
$ws.Range("E7").Value = 0.1737670030076259
$ws.Range("G7").Value = 0.00257062716290884
$ws.Range("J7").Value = 0.06839744146114146
$ws.Range("K7").Value = 0.5615798828134757
$ws.Range("L7").Value = 0.3795006175838651
$ws.Range("M7").Value = 0.3250902509863565
$ws.Range("N7").Value = 3.212370202317373
$ws.Range("O7").Value = 7.652756363687018
# Row 8
$ws.Range("B8").Value = 1.142357676429015
$ws.Range("C8").Value = 0.09105745775769236
$ws.Range("D8").Value = 0.5252852940892865
$ws.Range("E8").Value = 0.1730564314489769
$ws.Range("G8").Value = 0.002566861607763858
$ws.Range("J8").Value = 0.06834306735797036
$ws.Range("K8").Value = 0.5983275723927193
$ws.Range("L8").Value = 0.3812051550222932
$ws.Range("M8").Value = 0.3316218634536163
$ws.Range("N8").Value = 3.18091202994222
$ws.Range("O8").Value = 7.613199100485559
# Row 9
$ws.Range("B9").Value = 1.220796232579886
$ws.Range("C9").Value = 0.09445880840732457
$ws.Range("D9").Value = 0.5297797430536235
$ws.Range("E9").Value = 0.1720455726569465
$ws.Range("G9").Value = 0.002560233223262081
$ws.Range("J9").Value = 0.06825243302527895
$ws.Range("K9").Value = 0.672451865374228
$ws.Range("L9").Value = 0.3856466157283478
$ws.Range("M9").Value = 0.3455457261583348
$ws.Range("N9").Value = 3.125450385904813
$ws.Range("O9").Value = 7.558018571613218
# Row 10
$ws.Range("B10").Value = 1.28085783879186
$ws.Range("C10").Value = 0.09689327612632326
$ws.Range("D10").Value = 0.5340567364804656
$ws.Range("E10").Value = 0.1715346135220575
$ws.Range("G10").Value = 0.002555819565331291
$ws.Range("J10").Value = 0.06819550699499555
$ws.Range("K10").Value = 0.7281532288558878
$ws.Range("L10").Value = 0.3895750416492092
$ws.Range("M10").Value = 0.356449883292413
$ws.Range("N10").Value = 3.088490456556464
$ws.Range("O10").Value = 7.531092480789709
# Row 11
$ws.Range("B11").Value = 1.308705097279045
$ws.Range("C11").Value = 0.0979867987382832
$ws.Range("D11").Value = 0.536213046569344
$ws.Range("E11").Value = 0.1713521683017589
$ws.Range("G11").Value = 0.002553909768868521
$ws.Range("J11").Value = 0.0681716899144913
$ws.Range("K11").Value = 0.7537601734641441
$ws.Range("L11").Value = 0.3915057172687568
$ws.Range("M11").Value = 0.3615557709717834
$ws.Range("N11").Value = 3.072497020894247
$ws.Range("O11").Value = 7.521792491190553
# Row 12
$ws.Range("B12").Value = 1.319325072928933
$ws.Range("C12").Value = 0.09839888145141629
$ws.Range("D12").Value = 0.5370597713825447
$ws.Range("E12").Value = 0.1712902432977081
$ws.Range("G12").Value = 0.002553200596995375
$ws.Range("J12").Value = 0.06816296850045589
$ws.Range("K12").Value = 0.7634950392958331
$ws.Range("L12").Value = 0.3922573737982731
$ws.Range("M12").Value = 0.3635100415447283
$ws.Range("N12").Value = 3.066558512375906
$ws.Range("O12").Value = 7.518694227018727
# Row 13
$ws.Range("B13").Value = 1.317034553903085
$ws.Range("C13").Value = 0.09831022160802405
$ws.Range("D13").Value = 0.5368760738107454
$ws.Range("E13").Value = 0.1713032617874539
$ws.Range("G13").Value = 0.002553352707169651
$ws.Range("J13").Value = 0.06816483359851055
$ws.Range("K13").Value = 0.7613967762037817
$ws.Range("L13").Value = 0.3920945786704806
$ws.Range("M13").Value = 0.3630882324547002
$ws.Range("N13").Value = 3.06783223510606
$ws.Range("O13").Value = 7.519342670493927
# Row 14
$ws.Range("B14").Value = 1.309577313286724
$ws.Range("C14").Value = 0.09802074136919003
$ws.Range("D14").Value = 0.5362821030199001
$ws.Range("E14").Value = 0.1713469302839563
$ws.Range("G14").Value = 0.002553851144189111
$ws.Range("J14").Value = 0.06817096644172516
$ws.Range("K14").Value = 0.7545603064885711
$ws.Range("L14").Value = 0.391567145191857
$ws.Range("M14").Value = 0.3617161343247872
$ws.Range("N14").Value = 3.072006094330007
$ws.Range("O14").Value = 7.521529112168594
# Row 15
$ws.Range("B15").Value = 1.305019262651399
$ws.Range("C15").Value = 0.09784316430358331
$ws.Range("D15").Value = 0.5359222051597214
$ws.Range("E15").Value = 0.1713746105732987
$ws.Range("G15").Value = 0.002554158276061158
$ws.Range("J15").Value = 0.0681747617018722
$ws.Range("K15").Value = 0.7503777151071063
$ws.Range("L15").Value = 0.3912467504121508
$ws.Range("M15").Value = 0.3608783866765677
$ws.Range("N15").Value = 3.074578051415678
$ws.Range("O15").Value = 7.522923497571242
# Row 16
$ws.Range("B16").Value = 1.279048461396684
$ws.Range("C16").Value = 0.09682153091623746
$ws.Range("D16").Value = 0.5339200458623594
$ws.Range("E16").Value = 0.1715475405146911
$ws.Range("G16").Value = 0.002555946341560677
$ws.Range("J16").Value = 0.06819710521445899
$ws.Range("K16").Value = 0.7264851086803503
$ws.Range("L16").Value = 0.3894517487971143
$ws.Range("M16").Value = 0.3561191175501435
$ws.Range("N16").Value = 3.089552159376723
$ws.Range("O16").Value = 7.53175954370127
# Row 17
$ws.Range("B17").Value = 1.26325019312776
$ws.Range("C17").Value = 0.09619121897583227
$ws.Range("D17").Value = 0.5327456624519726
$ws.Range("E17").Value = 0.1716664133281736
$ws.Range("G17").Value = 0.002557068314905764
$ws.Range("J17").Value = 0.06821134372359605
$ws.Range("K17").Value = 0.7118960955945681
$ws.Range("L17").Value = 0.3883872854012225
$ws.Range("M17").Value = 0.353236630673301
$ws.Range("N17").Value = 3.098948242502175
$ws.Range("O17").Value = 7.537935010006635
# Row 18
$ws.Range("B18").Value = 1.254212895727164
$ws.Range("C18").Value = 0.09582736961399974
$ws.Range("D18").Value = 0.5320900259000894
$ws.Range("E18").Value = 0.1717394929099534
$ws.Range("G18").Value = 0.002557722872323575
$ws.Range("J18").Value = 0.06821972908834262
$ws.Range("K18").Value = 0.7035301462596237
$ws.Range("L18").Value = 0.3877885579000662
$ws.Range("M18").Value = 0.3515924073076064
$ws.Range("N18").Value = 3.104429788093917
$ws.Range("O18").Value = 7.541764567556385
# Row 19
$ws.Range("B19").Value = 1.251161535124169
$ws.Range("C19").Value = 0.09570395160702816
$ws.Range("D19").Value = 0.5318714494088397
$ws.Range("E19").Value = 0.1717650456357731
$ws.Range("G19").Value = 0.002557946081166067
$ws.Range("J19").Value = 0.0682226018904224
$ws.Range("K19").Value = 0.7007019335860036
$ws.Range("L19").Value = 0.3875881648157105
$ws.Range("M19").Value = 0.3510380601302714
$ws.Range("N19").Value = 3.106298999772889
$ws.Range("O19").Value = 7.543108885039544
# Row 20
$ws.Range("B20").Value = 1.264926832796135
$ws.Range("C20").Value = 0.09625845245556519
$ws.Range("D20").Value = 0.5328686253061932
$ws.Range("E20").Value = 0.1716532721364263
$ws.Range("G20").Value = 0.002556947924366774
$ws.Range("J20").Value = 0.06820980776010277
$ws.Range("K20").Value = 0.7134465095956273
$ws.Range("L20").Value = 0.3884992004219185
$ws.Range("M20").Value = 0.3535420587700813
$ws.Range("N20").Value = 3.097940027137867
$ws.Range("O20").Value = 7.537248895607263
# Row 21
$ws.Range("B21").Value = 1.311765660210796
$ws.Range("C21").Value = 0.09810582335928331
$ws.Range("D21").Value = 0.5364557485938519
$ws.Range("E21").Value = 0.1713339095905404
$ws.Range("G21").Value = 0.00255370436093456
$ws.Range("J21").Value = 0.06816915701089954
$ws.Range("K21").Value = 0.7565673139584419
$ws.Range("L21").Value = 0.3917215082227159
$ws.Range("M21").Value = 0.362118589955422
$ws.Range("N21").Value = 3.070776932283948
$ws.Range("O21").Value = 7.520875413944111
# Row 22
$ws.Range("B22").Value = 1.342813361845117
$ws.Range("C22").Value = 0.09930146437719145
$ws.Range("D22").Value = 0.5389759803592113
$ws.Range("E22").Value = 0.171166927478275
$ws.Range("G22").Value = 0.002551666233496622
$ws.Range("J22").Value = 0.06814432341194987
$ws.Range("K22").Value = 0.7849709388691508
$ws.Range("L22").Value = 0.3939472166581908
$ws.Range("M22").Value = 0.367844913182175
$ws.Range("N22").Value = 3.053711237005892
$ws.Range("O22").Value = 7.51264225591072
# Row 23
$ws.Range("B23").Value = 1.326202960057515
$ws.Range("C23").Value = 0.09866440340150717
$ws.Range("D23").Value = 0.5376148352602286
$ws.Range("E23").Value = 0.1712522382913804
$ws.Range("G23").Value = 0.002552746563829828
$ws.Range("J23").Value = 0.06815741934746367
$ws.Range("K23").Value = 0.7697912649645389
$ws.Range("L23").Value = 0.3927483906029607
$ws.Range("M23").Value = 0.3647776373842291
$ws.Range("N23").Value = 3.062756672409918
$ws.Range("O23").Value = 7.516810834236082
# Row 24
$ws.Range("B24").Value = 1.264168682730229
$ws.Range("C24").Value = 0.09622806079321577
$ws.Range("D24").Value = 0.5328129729512909
$ws.Range("E24").Value = 0.1716591985090457
$ws.Range("G24").Value = 0.002557002323453897
$ws.Range("J24").Value = 0.06821050154811648
$ws.Range("K24").Value = 0.7127455005173715
$ws.Range("L24").Value = 0.3884485623719343
$ws.Range("M24").Value = 0.3534039343578073
$ws.Range("N24").Value = 3.098395593306435
$ws.Range("O24").Value = 7.537558218178702
# Row 25
$ws.Range("B25").Value = 1.199147588733922
$ws.Range("C25").Value = 0.09354998561101979
$ws.Range("D25").Value = 0.5283922732091924
$ws.Range("E25").Value = 0.1722782361266635
$ws.Range("G25").Value = 0.002561945931828739
$ws.Range("J25").Value = 0.06827524886432812
$ws.Range("K25").Value = 0.6521800514976235
$ws.Range("L25").Value = 0.3843279195306835
$ws.Range("M25").Value = 0.3416601296652075
$ws.Range("N25").Value = 3.139788437071736
$ws.Range("O25").Value = 7.570552898909483
